$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updates to existing rows (AgTests / AgPosit columns F & G) ---
$ws.Range("F303").Value = 9330
$ws.Range("G303").Value = 590
$ws.Range("F306").Value = 77294
$ws.Range("G306").Value = 7709
$ws.Range("F307").Value = 75200
$ws.Range("G307").Value = 6314
$ws.Range("F308").Value = 15691
$ws.Range("G308").Value = 1028
$ws.Range("F309").Value = 78080
$ws.Range("G309").Value = 5533
$ws.Range("F310").Value = 79478
$ws.Range("G310").Value = 4050
$ws.Range("F311").Value = 61230
$ws.Range("G311").Value = 1918
$ws.Range("F313").Value = 76961
$ws.Range("G313").Value = 3468
$ws.Range("F314").Value = 65297
$ws.Range("F315").Value = 56781
$ws.Range("G315").Value = 2630
$ws.Range("F316").Value = 50671
$ws.Range("G316").Value = 2293
$ws.Range("F317").Value = 64145
$ws.Range("F318").Value = 49199
$ws.Range("G318").Value = 1128
$ws.Range("F320").Value = 74647
$ws.Range("G320").Value = 3385
$ws.Range("F321").Value = 91107
$ws.Range("G321").Value = 2670
$ws.Range("F322").Value = 110271
$ws.Range("G322").Value = 2328
$ws.Range("F323").Value = 217380
$ws.Range("G323").Value = 3103
$ws.Range("F324").Value = 248856
$ws.Range("G324").Value = 2848
$ws.Range("F325").Value = 781689
$ws.Range("G325").Value = 6515
$ws.Range("F326").Value = 407052
$ws.Range("G326").Value = 3793
$ws.Range("F327").Value = 225204
$ws.Range("G327").Value = 2738
$ws.Range("F328").Value = 181025
$ws.Range("G328").Value = 2667
$ws.Range("F329").Value = 73265
$ws.Range("G329").Value = 1727
$ws.Range("F330").Value = 72243
$ws.Range("G330").Value = 2082
$ws.Range("F331").Value = 154333
$ws.Range("G331").Value = 2716
$ws.Range("F332").Value = 486302
$ws.Range("G332").Value = 4819
$ws.Range("F333").Value = 256495
$ws.Range("G333").Value = 2865
$ws.Range("F334").Value = 192920
$ws.Range("F335").Value = 150566
$ws.Range("G335").Value = 3749
$ws.Range("F336").Value = 82021
$ws.Range("F338").Value = 221346
$ws.Range("G338").Value = 3048
$ws.Range("F340").Value = 388616
$ws.Range("G340").Value = 3315
$ws.Range("F341").Value = 286115
$ws.Range("G341").Value = 3631
$ws.Range("F342").Value = 178906
$ws.Range("F343").Value = 133973
$ws.Range("F344").Value = 136294
$ws.Range("F345").Value = 293817
$ws.Range("G345").Value = 3340
$ws.Range("F346").Value = 675792
$ws.Range("G346").Value = 4833
$ws.Range("F347").Value = 345754
$ws.Range("F348").Value = 234177
$ws.Range("G348").Value = 3283
$ws.Range("F349").Value = 161083
$ws.Range("G349").Value = 2761
$ws.Range("F350").Value = 128028
$ws.Range("F351").Value = 150259
$ws.Range("G351").Value = 2797
$ws.Range("F352").Value = 307939
$ws.Range("G352").Value = 3567
$ws.Range("F353").Value = 725556
$ws.Range("G353").Value = 5296
$ws.Range("F354").Value = 317006
$ws.Range("G354").Value = 2890
$ws.Range("F355").Value = 222292
$ws.Range("G355").Value = 3437
$ws.Range("F359").Value = 320690
$ws.Range("F362").Value = 228973
$ws.Range("F363").Value = 189530
$ws.Range("F364").Value = 168748
$ws.Range("G364").Value = 2490
$ws.Range("F366").Value = 338805
$ws.Range("F367").Value = 767372
$ws.Range("F372").Value = 178683
$ws.Range("G372").Value = 1859
$ws.Range("F373").Value = 350526
$ws.Range("F375").Value = 350333
$ws.Range("G375").Value = 1845
$ws.Range("F376").Value = 221795
$ws.Range("F378").Value = 157501
$ws.Range("F379").Value = 180801
$ws.Range("F380").Value = 345347
$ws.Range("G380").Value = 2028
$ws.Range("F381").Value = 748090
$ws.Range("F383").Value = 222648
$ws.Range("G383").Value = 1769
$ws.Range("F384").Value = 172350
$ws.Range("F387").Value = 351627
$ws.Range("F389").Value = 353571
$ws.Range("F393").Value = 308473
$ws.Range("F395").Value = 753096
$ws.Range("F401").Value = 272706
$ws.Range("F404").Value = 224153
$ws.Range("F406").Value = 171511
$ws.Range("F407").Value = 158330
$ws.Range("F408").Value = 304893
$ws.Range("F411").Value = 225429
$ws.Range("F413").Value = 149792
$ws.Range("F418").Value = 202443
$ws.Range("F421").Value = 153261
$ws.Range("F422").Value = 298536
$ws.Range("F429").Value = 171228
$ws.Range("F432").Value = 118199
$ws.Range("F438").Value = 118394
$ws.Range("F452").Value = 72357
$ws.Range("F455").Value = 50002
$ws.Range("F456").Value = 48179
$ws.Range("F458").Value = 67667
$ws.Range("F460").Value = 55983
$ws.Range("F461").Value = 43733
$ws.Range("F462").Value = 41965
$ws.Range("F463").Value = 44698
$ws.Range("F466").Value = 49379
$ws.Range("G466").Value = 55
$ws.Range("F470").Value = 41613
$ws.Range("F471").Value = 63292
$ws.Range("F472").Value = 48550
$ws.Range("F473").Value = 39023
$ws.Range("F474").Value = 43312
$ws.Range("F475").Value = 34995
$ws.Range("G475").Value = 27
$ws.Range("F476").Value = 35732
$ws.Range("F477").Value = 36298
$ws.Range("F478").Value = 52336
$ws.Range("F479").Value = 40545
$ws.Range("F480").Value = 33097
$ws.Range("F481").Value = 41985
$ws.Range("G481").Value = 35
$ws.Range("F482").Value = 35165
$ws.Range("F483").Value = 63251
$ws.Range("F484").Value = 8119
$ws.Range("F485").Value = 13748
$ws.Range("F486").Value = 8839
$ws.Range("F487").Value = 6786
$ws.Range("F488").Value = 6274
$ws.Range("F489").Value = 11885
$ws.Range("F490").Value = 10135
$ws.Range("G490").Value = 77
$ws.Range("F491").Value = 9405
$ws.Range("F492").Value = 12979
$ws.Range("F493").Value = 7566
$ws.Range("G493").Value = 7
$ws.Range("F494").Value = 6042
$ws.Range("G494").Value = 6

# --- New rows 495-497 (data through 2021-07-14) ---
$ws.Range("A495").Value = 44389
$ws.Range("A495").NumberFormat = "yyyy-mm-dd"
$ws.Range("B495").Value = 391953
$ws.Range("C495").Value = 6528
$ws.Range("D495").Value = 28
$ws.Range("E495").Value = 12522
$ws.Range("F495").Value = 9586
$ws.Range("G495").Value = 13
$ws.Range("A496").Value = 44390
$ws.Range("A496").NumberFormat = "yyyy-mm-dd"
$ws.Range("B496").Value = 391971
$ws.Range("C496").Value = 4697
$ws.Range("D496").Value = 18
$ws.Range("E496").Value = 12523
$ws.Range("F496").Value = 7435
$ws.Range("G496").Value = 14
$ws.Range("A497").Value = 44391
$ws.Range("A497").NumberFormat = "yyyy-mm-dd"
$ws.Range("B497").Value = 392000
$ws.Range("C497").Value = 5909
$ws.Range("D497").Value = 29
$ws.Range("E497").Value = 12524
$ws.Range("F497").Value = 4943
$ws.Range("G497").Value = 8
